
# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

function Get-ParaRange($d, $idx) {
    $p = $d.Paragraphs($idx)
    $start = $p.Range.Start
    $end = $p.Range.End - 1   # exclude the trailing paragraph mark
    return $d.Range($start, $end)
}

# Replace the *entire* text of paragraph $idx (paragraph mark preserved).
function Set-ParaText($d, $idx, $text) {
    $r = Get-ParaRange $d $idx
    $r.Text = $text
}

# Append $text right before the (preserved) paragraph mark of paragraph $idx.
function Add-ParaText($d, $idx, $text) {
    $p = $d.Paragraphs($idx)
    $pos = $p.Range.End - 1
    $r = $d.Range($pos, $pos)
    $r.InsertBefore($text)
}

# Append a manual line-break right before the paragraph mark, then return.
function Add-ParaBreak($d, $idx) {
    $p = $d.Paragraphs($idx)
    $pos = $p.Range.End - 1
    $r = $d.Range($pos, $pos)
    $r.InsertBreak(6) | Out-Null   # wdLineBreak
}

function Set-ParaStyle($d, $idx, $styleName) {
    $p = $d.Paragraphs($idx)
    $p.Style = $styleName
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title: "Big Fish" -> "Title"
# ---------------------------------------------------------------------------
Set-ParaText $d 1 "Title"

# ---------------------------------------------------------------------------
# 2. Author: "John August" -> "Author"
# ---------------------------------------------------------------------------
Set-ParaText $d 2 "Author"

# ---------------------------------------------------------------------------
# 3. Action: "This is a Southern story..." -> "Blah blah" + "."
# ---------------------------------------------------------------------------
Set-ParaText $d 3 "Blah blah"
Add-ParaText $d 3 "."

# ---------------------------------------------------------------------------
# 4. Action "====" is dropped entirely -- handled later (deleted last,
#    descending index order, after every other edit below).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 5. Transition "FADE IN:" is unchanged - nothing to do.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 6. Action: "A RIVER." -> "A HOUSE IN MAINE" + "."
# ---------------------------------------------------------------------------
Set-ParaText $d 6 "A HOUSE IN MAINE"
Add-ParaText $d 6 "."

# ---------------------------------------------------------------------------
# 7. Action: "We're underwater..." -> "Action description."
# ---------------------------------------------------------------------------
Set-ParaText $d 7 "Action description."

# ---------------------------------------------------------------------------
# 8. Action "This is The Beast." -> Character "CHARACTER" + " (V.O.)"
# ---------------------------------------------------------------------------
Set-ParaStyle $d 8 "Character"
Set-ParaText $d 8 "CHARACTER"
Add-ParaText $d 8 " (V.O.)"

# ---------------------------------------------------------------------------
# 9. Character "EDWARD (V.O.)" -> Dialogue "Scoobajy" + " " + "tralala" + "."
# ---------------------------------------------------------------------------
Set-ParaStyle $d 9 "Dialogue"
Set-ParaText $d 9 "Scoobajy"
Add-ParaText $d 9 " "
Add-ParaText $d 9 "tralala"
Add-ParaText $d 9 "."

# ---------------------------------------------------------------------------
# 10. Dialogue (2 runs) -> Parenthetical "(sighs)"
# ---------------------------------------------------------------------------
Set-ParaStyle $d 10 "Parenthetical"
Set-ParaText $d 10 "(sighs)"

# ---------------------------------------------------------------------------
# 11. Parenthetical "(sighs)" -> Dialogue "Nothing to say."
# ---------------------------------------------------------------------------
Set-ParaStyle $d 11 "Dialogue"
Set-ParaText $d 11 "Nothing to say."

# ---------------------------------------------------------------------------
# 12-15 dropped entirely -- handled later (deleted, descending order).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 16. Lyrics: "Some said ... <br/> a thief ..." ->
#     "I think I'll sing" <br/> "about a ring" + "."
# ---------------------------------------------------------------------------
$rightSingleQuote = [char]0x2019
$lyricsLine1 = "I think I" + $rightSingleQuote + "ll sing"
Set-ParaText $d 16 $lyricsLine1
Add-ParaBreak $d 16
Add-ParaText $d 16 "about a ring"
Add-ParaText $d 16 "."

# ---------------------------------------------------------------------------
# 17 dropped entirely -- handled below.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 18. SceneHeading (empty) is unchanged - nothing to do.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Now delete the surplus paragraphs, strictly from the highest original
# index to the lowest, so none of the earlier (already-edited) paragraph
# indices shift while we still need them.
# ---------------------------------------------------------------------------
$toDelete = @(17, 15, 14, 13, 12, 4)
foreach ($idx in $toDelete) {
    $d.Paragraphs($idx).Range.Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# Add the new "Note" character style (styles.xml addition).
# ---------------------------------------------------------------------------
$note = $d.Styles.Add("Note", 2)   # wdStyleTypeCharacter
$note.BaseStyle = "DefaultParagraphFont"
$note.Font.Hidden = $true
$note.Priority = 1

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
